$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value2 = 2199.077
$ws.Cells.Item(19, 9).Value2 = 3309.7144
$ws.Cells.Item(19, 11).Value2 = 3309.7144
$ws.Cells.Item(19, 13).Value2 = -3134.7144

$ws.Cells.Item(32, 8).Value2 = 641.1429000000001
$ws.Cells.Item(32, 9).Value2 = 399
$ws.Cells.Item(32, 10).Value2 = 738
$ws.Cells.Item(32, 11).Value2 = 399
$ws.Cells.Item(32, 12).Value2 = 738
$ws.Cells.Item(32, 13).Value2 = -73
$ws.Cells.Item(32, 14).Value2 = -1390

$ws.Cells.Item(129, 8).Value2 = 847.1884
$ws.Cells.Item(129, 10).Value2 = 849.48486
$ws.Cells.Item(129, 12).Value2 = 2548.45458
$ws.Cells.Item(129, 14).Value2 = -12548.45458

$ws.Cells.Item(132, 8).Value2 = 4067.1052
$ws.Cells.Item(132, 9).Value2 = 4485.9375
$ws.Cells.Item(132, 10).Value2 = 1833.3334
$ws.Cells.Item(132, 11).Value2 = 13457.8125
$ws.Cells.Item(132, 12).Value2 = 5500.0002
$ws.Cells.Item(132, 13).Value2 = -10927.8125
$ws.Cells.Item(132, 14).Value2 = -10560.0002

$ws.Cells.Item(137, 8).Value2 = 37895.355
$ws.Cells.Item(137, 9).Value2 = 2278.7646
$ws.Cells.Item(137, 10).Value2 = 92939.17999999999
$ws.Cells.Item(137, 11).Value2 = 6836.293799999999
$ws.Cells.Item(137, 12).Value2 = 278817.54
$ws.Cells.Item(137, 13).Value2 = -4286.293799999999
$ws.Cells.Item(137, 14).Value2 = -283917.54

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 23425.723
$ws.Cells.Item(32, 9).Value2 = 25286.605
$ws.Cells.Item(32, 11).Value2 = 25286.605
$ws.Cells.Item(32, 13).Value2 = -24999.605

$ws.Cells.Item(45, 8).Value2 = 3316
$ws.Cells.Item(45, 9).Value2 = 3204.625
$ws.Cells.Item(45, 11).Value2 = 3204.625
$ws.Cells.Item(45, 13).Value2 = -2827.625

$ws.Cells.Item(61, 8).Value2 = 4750.0713
$ws.Cells.Item(61, 9).Value2 = 4012.625
$ws.Cells.Item(61, 11).Value2 = 4012.625
$ws.Cells.Item(61, 13).Value2 = -3800.625

$ws.Cells.Item(74, 8).Value2 = 1716.4445
$ws.Cells.Item(74, 9).Value2 = 1946.5
$ws.Cells.Item(74, 10).Value2 = 1059.1428
$ws.Cells.Item(74, 11).Value2 = 1946.5
$ws.Cells.Item(74, 12).Value2 = 1059.1428
$ws.Cells.Item(74, 13).Value2 = -1072.5
$ws.Cells.Item(74, 14).Value2 = -2807.1428

$ws.Cells.Item(77, 8).Value2 = 1716.4445
$ws.Cells.Item(77, 9).Value2 = 1946.5
$ws.Cells.Item(77, 10).Value2 = 1059.1428
$ws.Cells.Item(77, 11).Value2 = 9732.5
$ws.Cells.Item(77, 12).Value2 = 5295.714
$ws.Cells.Item(77, 13).Value2 = -5364.5
$ws.Cells.Item(77, 14).Value2 = -14031.714

$ws.Cells.Item(110, 8).Value2 = 2660.8667
$ws.Cells.Item(110, 9).Value2 = 2954.4285
$ws.Cells.Item(110, 10).Value2 = 2404
$ws.Cells.Item(110, 11).Value2 = 2954.4285
$ws.Cells.Item(110, 12).Value2 = 2404
$ws.Cells.Item(110, 13).Value2 = -909.4285
$ws.Cells.Item(110, 14).Value2 = -6494

$ws.Cells.Item(132, 8).Value2 = 20545.607
$ws.Cells.Item(132, 9).Value2 = 2463.1538
$ws.Cells.Item(132, 10).Value2 = 36217.066
$ws.Cells.Item(132, 11).Value2 = 7389.4614
$ws.Cells.Item(132, 12).Value2 = 108651.198
$ws.Cells.Item(132, 13).Value2 = -4859.4614
$ws.Cells.Item(132, 14).Value2 = -113711.198

$ws.Cells.Item(136, 8).Value2 = 4750.0713
$ws.Cells.Item(136, 9).Value2 = 4012.625
$ws.Cells.Item(136, 11).Value2 = 12037.875
$ws.Cells.Item(136, 13).Value2 = -9487.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value2 = 2059.7273
$ws.Cells.Item(99, 9).Value2 = 1853.2222
$ws.Cells.Item(99, 10).Value2 = 2989
$ws.Cells.Item(99, 11).Value2 = 1853.2222
$ws.Cells.Item(99, 12).Value2 = 2989
$ws.Cells.Item(99, 13).Value2 = -355.2221999999999
$ws.Cells.Item(99, 14).Value2 = -5985

$ws.Cells.Item(107, 8).Value2 = 4304.3335
$ws.Cells.Item(107, 9).Value2 = 5000
$ws.Cells.Item(107, 10).Value2 = 3956.5
$ws.Cells.Item(107, 11).Value2 = 5000
$ws.Cells.Item(107, 12).Value2 = 3956.5
$ws.Cells.Item(107, 13).Value2 = -3080
$ws.Cells.Item(107, 14).Value2 = -7796.5

$ws.Cells.Item(111, 8).Value2 = 31900
$ws.Cells.Item(111, 10).Value2 = 31900
$ws.Cells.Item(111, 12).Value2 = 31900
$ws.Cells.Item(111, 14).Value2 = -40080

$ws.Cells.Item(134, 8).Value2 = 42002
$ws.Cells.Item(134, 9).Value2 = 49090.727
$ws.Cells.Item(134, 10).Value2 = 3014
$ws.Cells.Item(134, 11).Value2 = 147272.181
$ws.Cells.Item(134, 12).Value2 = 9042
$ws.Cells.Item(134, 13).Value2 = -144737.181
$ws.Cells.Item(134, 14).Value2 = -14112

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 10282.162
$ws.Cells.Item(31, 9).Value2 = 12556.821
$ws.Cells.Item(31, 10).Value2 = 3205.4443
$ws.Cells.Item(31, 11).Value2 = 12556.821
$ws.Cells.Item(31, 12).Value2 = 3205.4443
$ws.Cells.Item(31, 13).Value2 = -12261.821
$ws.Cells.Item(31, 14).Value2 = -3795.4443

$ws.Cells.Item(34, 8).Value2 = 10282.162
$ws.Cells.Item(34, 9).Value2 = 12556.821
$ws.Cells.Item(34, 10).Value2 = 3205.4443
$ws.Cells.Item(34, 11).Value2 = 12556.821
$ws.Cells.Item(34, 12).Value2 = 3205.4443
$ws.Cells.Item(34, 13).Value2 = -12354.821
$ws.Cells.Item(34, 14).Value2 = -3609.4443

$ws.Cells.Item(62, 8).Value2 = 6265.375
$ws.Cells.Item(62, 9).Value2 = 4368.3335
$ws.Cells.Item(62, 10).Value2 = 7403.6
$ws.Cells.Item(62, 11).Value2 = 4368.3335
$ws.Cells.Item(62, 12).Value2 = 7403.6
$ws.Cells.Item(62, 13).Value2 = -3744.3335
$ws.Cells.Item(62, 14).Value2 = -8651.6

$ws.Cells.Item(65, 8).Value2 = 6265.375
$ws.Cells.Item(65, 9).Value2 = 4368.3335
$ws.Cells.Item(65, 10).Value2 = 7403.6
$ws.Cells.Item(65, 11).Value2 = 21841.6675
$ws.Cells.Item(65, 12).Value2 = 37018
$ws.Cells.Item(65, 13).Value2 = -18721.6675
$ws.Cells.Item(65, 14).Value2 = -43258

$ws.Cells.Item(68, 8).Value2 = 69154
$ws.Cells.Item(68, 10).Value2 = 69154
$ws.Cells.Item(68, 12).Value2 = 69154
$ws.Cells.Item(68, 14).Value2 = -70652

$ws.Cells.Item(71, 8).Value2 = 69154
$ws.Cells.Item(71, 10).Value2 = 69154
$ws.Cells.Item(71, 12).Value2 = 207462
$ws.Cells.Item(71, 14).Value2 = -214950

$ws.Cells.Item(107, 8).Value2 = 2237.5715
$ws.Cells.Item(107, 9).Value2 = 1608.3334
$ws.Cells.Item(107, 11).Value2 = 1608.3334
$ws.Cells.Item(107, 13).Value2 = 311.6666

$ws.Cells.Item(134, 8).Value2 = 1055.0333
$ws.Cells.Item(134, 9).Value2 = 792.15
$ws.Cells.Item(134, 11).Value2 = 2376.45
$ws.Cells.Item(134, 13).Value2 = 158.5500000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value2 = 1325
$ws.Cells.Item(75, 9).Value2 = 1400
$ws.Cells.Item(75, 10).Value2 = 1250
$ws.Cells.Item(75, 11).Value2 = 4200
$ws.Cells.Item(75, 12).Value2 = 3750
$ws.Cells.Item(75, 13).Value2 = -3202
$ws.Cells.Item(75, 14).Value2 = -5746

$ws.Cells.Item(78, 8).Value2 = 1325
$ws.Cells.Item(78, 9).Value2 = 1400
$ws.Cells.Item(78, 10).Value2 = 1250
$ws.Cells.Item(78, 11).Value2 = 12600
$ws.Cells.Item(78, 12).Value2 = 11250
$ws.Cells.Item(78, 13).Value2 = -7608
$ws.Cells.Item(78, 14).Value2 = -21234

$ws.Cells.Item(131, 8).Value2 = 750.7
$ws.Cells.Item(131, 10).Value2 = 766.6667
$ws.Cells.Item(131, 12).Value2 = 2300.0001
$ws.Cells.Item(131, 14).Value2 = -12380.0001

$ws.Cells.Item(141, 8).Value2 = 2954.1667
$ws.Cells.Item(141, 9).Value2 = 2988
$ws.Cells.Item(141, 10).Value2 = 2911.875
$ws.Cells.Item(141, 11).Value2 = 8964
$ws.Cells.Item(141, 12).Value2 = 8735.625
$ws.Cells.Item(141, 13).Value2 = -3784
$ws.Cells.Item(141, 14).Value2 = -19095.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value2 = 3500
$ws.Cells.Item(113, 9).Value2 = 2350
$ws.Cells.Item(113, 10).Value2 = 5033.3335
$ws.Cells.Item(113, 11).Value2 = 2350
$ws.Cells.Item(113, 12).Value2 = 5033.3335
$ws.Cells.Item(113, 13).Value2 = -180
$ws.Cells.Item(113, 14).Value2 = -9373.333500000001

$ws.Cells.Item(126, 8).Value2 = 4721.1724
$ws.Cells.Item(126, 10).Value2 = 5760.933
$ws.Cells.Item(126, 12).Value2 = 17282.799
$ws.Cells.Item(126, 14).Value2 = -22222.799

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value2 = 2134.8333
$ws.Cells.Item(22, 9).Value2 = 3533.3333
$ws.Cells.Item(22, 10).Value2 = 1668.6666
$ws.Cells.Item(22, 11).Value2 = 3533.3333
$ws.Cells.Item(22, 12).Value2 = 1668.6666
$ws.Cells.Item(22, 13).Value2 = -3238.3333
$ws.Cells.Item(22, 14).Value2 = -2258.6666

$ws.Cells.Item(27, 8).Value2 = 2134.8333
$ws.Cells.Item(27, 9).Value2 = 3533.3333
$ws.Cells.Item(27, 10).Value2 = 1668.6666
$ws.Cells.Item(27, 11).Value2 = 3533.3333
$ws.Cells.Item(27, 12).Value2 = 1668.6666
$ws.Cells.Item(27, 13).Value2 = -3426.3333
$ws.Cells.Item(27, 14).Value2 = -1882.6666

$ws.Cells.Item(40, 8).Value2 = 5193.7856
$ws.Cells.Item(40, 9).Value2 = 5217.3335
$ws.Cells.Item(40, 11).Value2 = 5217.3335
$ws.Cells.Item(40, 13).Value2 = -5081.3335

$ws.Cells.Item(46, 8).Value2 = 2405.7144
$ws.Cells.Item(46, 9).Value2 = 2313
$ws.Cells.Item(46, 10).Value2 = 2475.25
$ws.Cells.Item(46, 11).Value2 = 2313
$ws.Cells.Item(46, 12).Value2 = 2475.25
$ws.Cells.Item(46, 13).Value2 = -2125
$ws.Cells.Item(46, 14).Value2 = -2851.25

$ws.Cells.Item(68, 8).Value2 = 4242
$ws.Cells.Item(68, 9).Value2 = 1962.6
$ws.Cells.Item(68, 10).Value2 = 6141.5
$ws.Cells.Item(68, 11).Value2 = 1962.6
$ws.Cells.Item(68, 12).Value2 = 6141.5
$ws.Cells.Item(68, 13).Value2 = -1213.6
$ws.Cells.Item(68, 14).Value2 = -7639.5

$ws.Cells.Item(71, 8).Value2 = 4242
$ws.Cells.Item(71, 9).Value2 = 1962.6
$ws.Cells.Item(71, 10).Value2 = 6141.5
$ws.Cells.Item(71, 11).Value2 = 9813
$ws.Cells.Item(71, 12).Value2 = 30707.5
$ws.Cells.Item(71, 13).Value2 = -6069
$ws.Cells.Item(71, 14).Value2 = -38195.5

$ws.Cells.Item(93, 8).Value2 = 1933.4736
$ws.Cells.Item(93, 9).Value2 = 1925.6471
$ws.Cells.Item(93, 11).Value2 = 1925.6471
$ws.Cells.Item(93, 13).Value2 = -677.6470999999999

$ws.Cells.Item(104, 8).Value2 = 22498
$ws.Cells.Item(104, 10).Value2 = 22498
$ws.Cells.Item(104, 12).Value2 = 22498
$ws.Cells.Item(104, 14).Value2 = -29486

$ws.Cells.Item(132, 8).Value2 = 1871.9048
$ws.Cells.Item(132, 9).Value2 = 1268.3572
$ws.Cells.Item(132, 11).Value2 = 3805.0716
$ws.Cells.Item(132, 13).Value2 = -1275.0716

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value2 = 1985.7142
$ws.Cells.Item(81, 9).Value2 = 1985.7142
$ws.Cells.Item(81, 10).Value2 = 0
$ws.Cells.Item(81, 11).Value2 = 3971.4284
$ws.Cells.Item(81, 12).Value2 = 0
$ws.Cells.Item(81, 13).Value2 = -2910.4284
$ws.Cells.Item(81, 14).Value2 = $null

$ws.Cells.Item(84, 8).Value2 = 1985.7142
$ws.Cells.Item(84, 9).Value2 = 1985.7142
$ws.Cells.Item(84, 10).Value2 = 0
$ws.Cells.Item(84, 11).Value2 = 19857.142
$ws.Cells.Item(84, 12).Value2 = 0
$ws.Cells.Item(84, 13).Value2 = -14553.142
$ws.Cells.Item(84, 14).Value2 = $null

$ws.Cells.Item(113, 8).Value2 = 3004232.5
$ws.Cells.Item(113, 10).Value2 = 9009169
$ws.Cells.Item(113, 12).Value2 = 27027507
$ws.Cells.Item(113, 14).Value2 = -27031847

$ws.Cells.Item(122, 8).Value2 = 1487.8
$ws.Cells.Item(122, 9).Value2 = 1464.0555
$ws.Cells.Item(122, 10).Value2 = 1548.8572
$ws.Cells.Item(122, 11).Value2 = 4392.166499999999
$ws.Cells.Item(122, 12).Value2 = 4646.571599999999
$ws.Cells.Item(122, 13).Value2 = -1942.166499999999
$ws.Cells.Item(122, 14).Value2 = -9546.571599999999

$ws.Cells.Item(126, 8).Value2 = 2681.3333
$ws.Cells.Item(126, 9).Value2 = 1200
$ws.Cells.Item(126, 11).Value2 = 3600
$ws.Cells.Item(126, 13).Value2 = -1130

Write-Host "done"
